# Fix mojibake in the Regional Economic Communities footnote (A103)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Update recalculated population-growth figures for rows 67-73 and 97-98
# Row 67
$ws.Range("F67").Value = 2.5921470390226
$ws.Range("G67").Value = 2.3831369160084899
$ws.Range("K67").Value = 2.5363805249015501
$ws.Range("L67").Value = 2.5903253279856702
$ws.Range("M67").Value = 2.6102984774830902
$ws.Range("N67").Value = 2.6353259164121301
$ws.Range("O67").Value = 2.6616111918730501
$ws.Range("R67").Value = 2.66855184353583
$ws.Range("S67").Value = 2.6859783177304499
$ws.Range("U67").Value = 2.7282096672888798
$ws.Range("W67").Value = 2.7235365249281198
$ws.Range("AA67").Value = 2.7268727227936398
$ws.Range("AB67").Value = 2.6783511273101901
$ws.Range("AC67").Value = 2.6896664392773499
$ws.Range("AE67").Value = 2.6405015333436999
$ws.Range("AH67").Value = 2.4872737504042499
$ws.Range("AI67").Value = 2.42115407360322
$ws.Range("AJ67").Value = 2.4081548816616198
$ws.Range("AL67").Value = 2.3812237511343
$ws.Range("AM67").Value = 2.35140904811455
$ws.Range("AN67").Value = 2.32382821362576
$ws.Range("AO67").Value = 2.29478300350006

# Row 68
$ws.Range("C68").Value = 2.5182118448480599
$ws.Range("D68").Value = 2.3494548731092499
$ws.Range("G68").Value = 2.4952531939674198
$ws.Range("H68").Value = 2.5187729724331902
$ws.Range("I68").Value = 2.4752302473691499
$ws.Range("J68").Value = 2.4301128067751199
$ws.Range("Q68").Value = 2.5646962251603602
$ws.Range("R68").Value = 2.5454256091826499
$ws.Range("U68").Value = 2.5506025382990201
$ws.Range("V68").Value = 2.5435590965137398
$ws.Range("W68").Value = 2.5498810990449901
$ws.Range("X68").Value = 0.50238290465549995
$ws.Range("Y68").Value = 2.41682304157733
$ws.Range("Z68").Value = 2.5355964455723501
$ws.Range("AA68").Value = 2.5563172282251401
$ws.Range("AC68").Value = 2.5127911234293698
$ws.Range("AH68").Value = 2.3194634493939099
$ws.Range("AI68").Value = 2.2741838866639399
$ws.Range("AK68").Value = 2.2665126918071299
$ws.Range("AL68").Value = 2.2490969809372499
$ws.Range("AN68").Value = 2.19577840803864

# Row 69
$ws.Range("C69").Value = 3.0863707907476101
$ws.Range("K69").Value = 2.6397075638872902
$ws.Range("L69").Value = 2.7847363745344502
$ws.Range("U69").Value = 3.0657085533914801
$ws.Range("AJ69").Value = 2.7825941515168502
$ws.Range("AK69").Value = 2.7739237059653599
$ws.Range("AM69").Value = 2.7244584258577098
$ws.Range("AN69").Value = 2.7014521015041
$ws.Range("AO69").Value = 2.6767858329292502

# Row 70
$ws.Range("E70").Value = 3.07321114707839
$ws.Range("F70").Value = 2.7091269490991099
$ws.Range("G70").Value = 2.1047578309432802
$ws.Range("H70").Value = 2.5400969063193801
$ws.Range("I70").Value = 3.3274598562155902
$ws.Range("J70").Value = 2.9609089969104398
$ws.Range("R70").Value = 3.2399172756611101
$ws.Range("S70").Value = 3.2349955225158098
$ws.Range("T70").Value = 3.2555804782918401
$ws.Range("W70").Value = 3.3814532545247298
$ws.Range("Y70").Value = 3.2785106176283598
$ws.Range("Z70").Value = 3.2829478422291798
$ws.Range("AB70").Value = 3.1802475126563299
$ws.Range("AC70").Value = 3.1913068392016801
$ws.Range("AD70").Value = 3.2128477144055099
$ws.Range("AF70").Value = 3.11966859741137
$ws.Range("AG70").Value = 3.0767709001535701
$ws.Range("AJ70").Value = 3.0053396862836901
$ws.Range("AM70").Value = 2.9377982728513001
$ws.Range("AN70").Value = 2.9094994619164498

# Row 72
$ws.Range("D72").Value = 2.8616576270849499
$ws.Range("E72").Value = 2.6848152531427498
$ws.Range("G72").Value = 3.1916274833956302
$ws.Range("H72").Value = 3.0676798615584699
$ws.Range("I72").Value = 2.9696240640787201
$ws.Range("J72").Value = 2.92163198639388
$ws.Range("V72").Value = 2.9239749896
$ws.Range("W72").Value = 2.8741537286101999
$ws.Range("AB72").Value = 2.63283686179554
$ws.Range("AC72").Value = 2.6681517907985999
$ws.Range("AD72").Value = 2.6245775755254299
$ws.Range("AF72").Value = 2.6732038833086902
$ws.Range("AG72").Value = 2.67703322663957
$ws.Range("AI72").Value = 2.52486399865648

# Row 73
$ws.Range("C73").Value = 2.9929293912671602
$ws.Range("F73").Value = 2.6003272255059202
$ws.Range("H73").Value = 2.8352038041092902
$ws.Range("K73").Value = 2.2328449392059699
$ws.Range("S73").Value = 2.4845318232602698
$ws.Range("U73").Value = 2.5501762760524298
$ws.Range("AB73").Value = 2.9662800154339801
$ws.Range("AC73").Value = 2.8173617961732198
$ws.Range("AE73").Value = 2.7435556095160898
$ws.Range("AF73").Value = 2.69169190775176
$ws.Range("AK73").Value = 2.58974461795842
$ws.Range("AL73").Value = 2.5803278867534898
$ws.Range("AM73").Value = 2.5493089973206402
$ws.Range("AN73").Value = 2.5220053291748399
$ws.Range("AP73").Value = 2.54708711239264

# Row 97
$ws.Range("C97").Value = 2.8712840883372999
$ws.Range("D97").Value = 2.9249276392543799
$ws.Range("E97").Value = 2.92669754493402
$ws.Range("F97").Value = 2.8679083605082099
$ws.Range("G97").Value = 2.9526136339493001
$ws.Range("H97").Value = 2.8688303578960501
$ws.Range("I97").Value = 2.7122396343194399
$ws.Range("J97").Value = 2.73054027080883
$ws.Range("K97").Value = 2.76854122895696
$ws.Range("L97").Value = 2.7630130356424401
$ws.Range("M97").Value = 2.7783243053547402
$ws.Range("N97").Value = 2.7984672466235199
$ws.Range("O97").Value = 2.8118953207824302
$ws.Range("P97").Value = 2.80383554361809
$ws.Range("Q97").Value = 2.7972054870211598
$ws.Range("R97").Value = 2.8074548437281699
$ws.Range("S97").Value = 2.8239773645694699
$ws.Range("T97").Value = 2.8382150346331101
$ws.Range("U97").Value = 2.8361417660871302
$ws.Range("V97").Value = 2.8198287773505002
$ws.Range("W97").Value = 2.8334399419298402
$ws.Range("X97").Value = 2.7920774270478099
$ws.Range("Y97").Value = 2.7794770725853901
$ws.Range("Z97").Value = 2.8195356430634102
$ws.Range("AA97").Value = 2.8140527585647699
$ws.Range("AB97").Value = 2.8111596535637098
$ws.Range("AC97").Value = 2.82615672421647
$ws.Range("AD97").Value = 2.8205769295637499
$ws.Range("AE97").Value = 2.7550839038070798
$ws.Range("AF97").Value = 2.7086479973325499
$ws.Range("AG97").Value = 2.6924563470530001
$ws.Range("AH97").Value = 2.6444346614751502
$ws.Range("AI97").Value = 2.5910214364820701
$ws.Range("AJ97").Value = 2.5738574864938402
$ws.Range("AK97").Value = 2.5606704324398799
$ws.Range("AL97").Value = 2.53350918494584
$ws.Range("AM97").Value = 2.5018843163832098
$ws.Range("AN97").Value = 2.4714644495370299
$ws.Range("AO97").Value = 2.44038948714815
$ws.Range("AP97").Value = 2.5015746351210599

# Row 98
$ws.Range("C98").Value = 2.7580747789735902
$ws.Range("D98").Value = 2.6697799860242202
$ws.Range("E98").Value = 2.2592175214822601
$ws.Range("F98").Value = 2.0139036975867599
$ws.Range("G98").Value = 2.0482906465511701
$ws.Range("H98").Value = 2.1411359890695798
$ws.Range("I98").Value = 2.1095529609313202
$ws.Range("J98").Value = 2.07856312472203
$ws.Range("K98").Value = 2.0831190212121
$ws.Range("L98").Value = 2.0718258363711799
$ws.Range("M98").Value = 2.13999660428088
$ws.Range("N98").Value = 2.1508015956938999
$ws.Range("O98").Value = 1.8489449554592901
$ws.Range("P98").Value = 1.7060571061198
$ws.Range("Q98").Value = 1.7815451200724
$ws.Range("R98").Value = 1.7092741296451901
$ws.Range("S98").Value = 1.59230617603761
$ws.Range("T98").Value = 1.54736693405895
$ws.Range("U98").Value = 1.4874924144758299
$ws.Range("V98").Value = 1.49593650226638
$ws.Range("W98").Value = 1.5903689934171601
$ws.Range("X98").Value = 1.55037376729301
$ws.Range("Y98").Value = 1.45425020245993
$ws.Range("Z98").Value = 1.3942330316978599
$ws.Range("AA98").Value = 1.3972830940858001
$ws.Range("AB98").Value = 1.3854609781049501
$ws.Range("AC98").Value = 1.2867654663596
$ws.Range("AD98").Value = 1.2123487023433499
$ws.Range("AE98").Value = 1.14869949335288
$ws.Range("AF98").Value = 1.11090440242374
$ws.Range("AG98").Value = 1.1834309887612999
$ws.Range("AH98").Value = 1.22828707688809
$ws.Range("AI98").Value = 1.28069165989329
$ws.Range("AJ98").Value = 1.3528737753707401
$ws.Range("AK98").Value = 1.34477818678824
$ws.Range("AL98").Value = 1.3128623776247601
$ws.Range("AM98").Value = 1.2793411274954201
$ws.Range("AN98").Value = 1.2470957979644699
$ws.Range("AO98").Value = 1.21512379955495
$ws.Range("AP98").Value = 1.2798298233897001
